# Update res_bus/vm_pu.xlsx values for the "380 kV" case (Case_3_151).
# Rows 2-25 correspond to bus indices 0-23; columns B:F and I:N hold the
# per-scenario voltage magnitudes (column G stays at 1, column H has no data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object "object[,]" 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.074117753921424
$dataBF[0,2] = 1.075303122240978
$dataBF[0,3] = 1.077924277549771
$dataBF[0,4] = 1.088132690176628
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.075356250915761
$dataBF[1,2] = 1.076293566355411
$dataBF[1,3] = 1.079025526235562
$dataBF[1,4] = 1.089297827551951
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.076157370759789
$dataBF[2,2] = 1.076934149521629
$dataBF[2,3] = 1.079738037495748
$dataBF[2,4] = 1.090051840923812
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.076494098824472
$dataBF[3,2] = 1.077203379757478
$dataBF[3,3] = 1.080037561626271
$dataBF[3,4] = 1.090368851336053
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.076550633329947
$dataBF[4,2] = 1.077248580568317
$dataBF[4,3] = 1.08008785216629
$dataBF[4,4] = 1.090422080175919
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.076161870380542
$dataBF[5,2] = 1.076937747266675
$dataBF[5,3] = 1.07974203981162
$dataBF[5,4] = 1.090056076739746
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.074536367025702
$dataBF[6,2] = 1.075637910264264
$dataBF[6,3] = 1.078296464986147
$dataBF[6,4] = 1.088526435608233
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.071669858693972
$dataBF[7,2] = 1.073345085876739
$dataBF[7,3] = 1.075748584162104
$dataBF[7,4] = 1.085831656811966
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.069757267812235
$dataBF[8,2] = 1.071814894170073
$dataBF[8,3] = 1.074049518608715
$dataBF[8,4] = 1.084035488413407
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.068928689085396
$dataBF[9,2] = 1.071151899492413
$dataBF[9,3] = 1.073313672547751
$dataBF[9,4] = 1.083257791717844
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.06862085315885
$dataBF[10,2] = 1.070905570499332
$dataBF[10,3] = 1.073040323915549
$dataBF[10,4] = 1.082968927212298
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.068686887989428
$dataBF[11,2] = 1.070958411789515
$dataBF[11,3] = 1.073098959176802
$dataBF[11,4] = 1.083030889355448
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.0689032446109
$dataBF[12,2] = 1.071131539159487
$dataBF[12,3] = 1.073291077936202
$dataBF[12,4] = 1.08323391396169
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.069036540387281
$dataBF[13,2] = 1.071238200222387
$dataBF[13,3] = 1.073409445593782
$dataBF[13,4] = 1.083359004954217
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.069812248898844
$dataBF[14,2] = 1.071858886151214
$dataBF[14,3] = 1.074098351205928
$dataBF[14,4] = 1.084087102613918
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.070298717570813
$dataBF[15,2] = 1.072248114847942
$dataBF[15,3] = 1.074530444965653
$dataBF[15,4] = 1.084543832817319
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.07058242672819
$dataBF[16,2] = 1.07247510574791
$dataBF[16,3] = 1.074782464604428
$dataBF[16,4] = 1.084810241347708
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.070679157550191
$dataBF[17,2] = 1.072552497101995
$dataBF[17,3] = 1.074868394576199
$dataBF[17,4] = 1.084901080762337
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.07024652822218
$dataBF[18,2] = 1.072206358401008
$dataBF[18,3] = 1.074484086787677
$dataBF[18,4] = 1.084494829430682
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.068839534791209
$dataBF[19,2] = 1.07108055922152
$dataBF[19,3] = 1.073234504388594
$dataBF[19,4] = 1.083174128105247
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.067954525396806
$dataBF[20,2] = 1.070372357945461
$dataBF[20,3] = 1.072448710502291
$dataBF[20,4] = 1.082343788508294
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.068423722225085
$dataBF[21,2] = 1.070747824069815
$dataBF[21,3] = 1.072865287770776
$dataBF[21,4] = 1.0827839642472
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.070270110465313
$dataBF[22,2] = 1.072225226464487
$dataBF[22,3] = 1.074505034094148
$dataBF[22,4] = 1.084516971934784
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.072411189654481
$dataBF[23,2] = 1.073938120586088
$dataBF[23,3] = 1.076407350903864
$dataBF[23,4] = 1.085831656811966
$ws.Range("B2:F25").Value = $dataBF

$dataIN = New-Object "object[,]" 24,6
$dataIN[0,0] = 1.061025189854597
$dataIN[0,1] = 1.079028844005987
$dataIN[0,2] = 1.077990048447023
$dataIN[0,3] = 1.080604299558871
$dataIN[0,4] = 1.090786165081943
$dataIN[0,5] = 1.080561188730813
$dataIN[1,0] = 1.061461218559055
$dataIN[1,1] = 1.079924930762111
$dataIN[1,2] = 1.078797357228423
$dataIN[1,3] = 1.081522641821931
$dataIN[1,4] = 1.091770165641073
$dataIN[1,5] = 1.081458548032914
$dataIN[2,0] = 1.061741916479581
$dataIN[2,1] = 1.080503925560297
$dataIN[2,2] = 1.079318820725571
$dataIN[2,3] = 1.082116202893224
$dataIN[2,4] = 1.092406383857647
$dataIN[2,5] = 1.082038365070125
$dataIN[3,0] = 1.061859577408084
$dataIN[3,1] = 1.080747136733662
$dataIN[3,2] = 1.079537824943205
$dataIN[3,3] = 1.082365577272592
$dataIN[3,4] = 1.09267373244588
$dataIN[3,5] = 1.082281921631253
$dataIN[4,0] = 1.061879313015596
$dataIN[4,1] = 1.080787961371102
$dataIN[4,2] = 1.079574583900776
$dataIN[4,3] = 1.082407439056352
$dataIN[4,4] = 1.092718614591351
$dataIN[4,5] = 1.08232280424436
$dataIN[5,0] = 1.06174349002336
$dataIN[5,1] = 1.080507176137293
$dataIN[5,2] = 1.079321747930014
$dataIN[5,3] = 1.082119535667981
$dataIN[5,4] = 1.092409956642984
$dataIN[5,5] = 1.082041620263313
$dataIN[6,0] = 1.061172846750539
$dataIN[6,1] = 1.079331853999309
$dataIN[6,2] = 1.078263073195009
$dataIN[6,3] = 1.08091479660095
$dataIN[6,4] = 1.091118816076781
$dataIN[6,5] = 1.080864629033082
$dataIN[7,0] = 1.060156226598859
$dataIN[7,1] = 1.077254348583314
$dataIN[7,2] = 1.076390465776289
$dataIN[7,3] = 1.078786725230455
$dataIN[7,4] = 1.08883980994373
$dataIN[7,5] = 1.078784173321103
$dataIN[8,0] = 1.059470989369555
$dataIN[8,1] = 1.075864936357153
$dataIN[8,2] = 1.075137228270776
$dataIN[8,3] = 1.077364460854618
$dataIN[8,4] = 1.087317806349881
$dataIN[8,5] = 1.077392787970213
$dataIN[9,0] = 1.059172486978762
$dataIN[9,1] = 1.075262241312739
$dataIN[9,2] = 1.074593400700009
$dataIN[9,3] = 1.076747745038367
$dataIN[9,4] = 1.086658111817526
$dataIN[9,5] = 1.076789237029707
$dataIN[10,0] = 1.059061339978127
$dataIN[10,1] = 1.075038211033826
$dataIN[10,2] = 1.074391222086414
$dataIN[10,3] = 1.076518537615505
$dataIN[10,4] = 1.086412971568899
$dataIN[10,5] = 1.0765648886021
$dataIN[11,0] = 1.059085193617219
$dataIN[11,1] = 1.075086273670011
$dataIN[11,2] = 1.07443459811091
$dataIN[11,3] = 1.076567709371776
$dataIN[11,4] = 1.086465559550569
$dataIN[11,5] = 1.076613019492741
$dataIN[12,0] = 1.059163305044967
$dataIN[12,1] = 1.075243726224833
$dataIN[12,2] = 1.074576692168569
$dataIN[12,3] = 1.076728801377441
$dataIN[12,4] = 1.086637850504167
$dataIN[12,5] = 1.076770695648251
$dataIN[13,0] = 1.059211396277241
$dataIN[13,1] = 1.075340716356384
$dataIN[13,2] = 1.074664217531003
$dataIN[13,3] = 1.076828037977877
$dataIN[13,4] = 1.086743991303972
$dataIN[13,5] = 1.076867823516914
$dataIN[14,0] = 1.059490762190796
$dataIN[14,1] = 1.075904912546291
$dataIN[14,2] = 1.075173295595514
$dataIN[14,3] = 1.077405371870239
$dataIN[14,4] = 1.087361574123501
$dataIN[14,5] = 1.077432820930125
$dataIN[15,0] = 1.059665521087581
$dataIN[15,1] = 1.076258530144972
$dataIN[15,2] = 1.075492313235667
$dataIN[15,3] = 1.077767285191132
$dataIN[15,4] = 1.087748790480335
$dataIN[15,5] = 1.077786940706356
$dataIN[16,0] = 1.059767282338338
$dataIN[16,1] = 1.076464686146669
$dataIN[16,2] = 1.0756782782833
$dataIN[16,3] = 1.077978299665804
$dataIN[16,4] = 1.087974583956989
$dataIN[16,5] = 1.077993389473222
$dataIN[17,0] = 1.059801951038621
$dataIN[17,1] = 1.076534962568546
$dataIN[17,2] = 1.075741668512472
$dataIN[17,3] = 1.07805023604353
$dataIN[17,4] = 1.088051563014097
$dataIN[17,5] = 1.07806376569568
$dataIN[18,0] = 1.059646788978177
$dataIN[18,1] = 1.07622060100682
$dataIN[18,2] = 1.075458097314319
$dataIN[18,3] = 1.077728463960315
$dataIN[18,4] = 1.087707252387217
$dataIN[18,5] = 1.077748957704478
$dataIN[19,0] = 1.05914031063535
$dataIN[19,1] = 1.075197364871359
$dataIN[19,2] = 1.074534853903433
$dataIN[19,3] = 1.07668136744565
$dataIN[19,4] = 1.086587117895177
$dataIN[19,5] = 1.076724268456338
$dataIN[20,0] = 1.058820305432354
$dataIN[20,1] = 1.074553074660535
$dataIN[20,2] = 1.073953350270393
$dataIN[20,3] = 1.07602225347877
$dataIN[20,4] = 1.085882263615009
$dataIN[20,5] = 1.076079063279513
$dataIN[21,0] = 1.058990094645563
$dataIN[21,1] = 1.074894714857349
$dataIN[21,2] = 1.07426171388102
$dataIN[21,3] = 1.076371735109154
$dataIN[21,4] = 1.086255975847904
$dataIN[21,5] = 1.076421188644593
$dataIN[22,0] = 1.059655253745611
$dataIN[22,1] = 1.076237739870507
$dataIN[22,2] = 1.075473558364169
$dataIN[22,3] = 1.077746005861364
$dataIN[22,4] = 1.087726021859085
$dataIN[22,5] = 1.077766120907318
$dataIN[23,0] = 1.060420364215353
$dataIN[23,1] = 1.077792204683998
$dataIN[23,2] = 1.076875426497012
$dataIN[23,3] = 1.079337502846399
$dataIN[23,4] = 1.089429451959823
$dataIN[23,5] = 1.079322793239143
$ws.Range("I2:N25").Value = $dataIN

Write-Host "vm_pu values updated for 380 kV case"